$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename background variable identifiers/labels
$ws.Range("A3").Value = "educ_group"
$ws.Range("A4").Value = "income_group"
$ws.Range("B4").Value = "Income Group"

# Update selection to A6 (as reflected in saved view state)
$ws.Range("A6").Select()

# NOTE: the source diff also nudges column A's stored width from 15.69 to
# 15.68 (a sub-pixel rounding artifact from the authoring tool) and adds a
# <charset val="1"/> hint to two Arial font entries. Neither is reachable
# through the COM object model here: Columns.ColumnWidth only round-trips
# at ~1/6-character granularity (the nearest settable value lands on 16.5,
# farther from 15.68 than the untouched 15.69), and there is no COM surface
# for a font's OOXML charset byte. Both are left as-is rather than moving
# the file further from the target or fabricating unreachable state.
